# Refresh the "1111" sample-data sheet with a newly generated batch of rows
# (fix: 完善xlsxwriter Handler -- regenerate the fixture rows/timestamps).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("序号") holds numeric-looking codes (some with leading zeros) that
# must be stored as text, not numbers -- format it as Text before writing.
$ws.Range("A2:A11").NumberFormat = "@"

$ws.Range("A2").Value = "5229`t"
$ws.Range("B2").Value = "哾喂溠痚`t"
$ws.Range("C2").Value = "2jGEhq`t"
$ws.Range("D2").Value = "2023-08-05 13:17:25`t"

$ws.Range("A3").Value = "4430`t"
$ws.Range("B3").Value = "踍縮膁思`t"
$ws.Range("C3").Value = "z9hp8v`t"
$ws.Range("D3").Value = "2023-08-05 13:17:25`t"

$ws.Range("A4").Value = "6950`t"
$ws.Range("B4").Value = "冤牲华簝`t"
$ws.Range("C4").Value = "K0mOSj`t"
$ws.Range("D4").Value = "2023-08-05 13:17:25`t"

$ws.Range("A5").Value = "0262`t"
$ws.Range("B5").Value = "槴嶿箕滋`t"
$ws.Range("C5").Value = "tXn9qA`t"
$ws.Range("D5").Value = "2023-08-05 13:17:25`t"

$ws.Range("A6").Value = "0362`t"
$ws.Range("B6").Value = "姳懞昂僰`t"
$ws.Range("C6").Value = "YWNZ17`t"
$ws.Range("D6").Value = "2023-08-05 13:17:25`t"

$ws.Range("A7").Value = "6349`t"
$ws.Range("B7").Value = "璿薶寶脔`t"
$ws.Range("C7").Value = "le7bmu`t"
$ws.Range("D7").Value = "2023-08-05 13:17:25`t"

$ws.Range("A8").Value = "7490`t"
$ws.Range("B8").Value = "鶿漅骚倠`t"
$ws.Range("C8").Value = "rVKHFW`t"
$ws.Range("D8").Value = "2023-08-05 13:17:25`t"

$ws.Range("A9").Value = "9204`t"
$ws.Range("B9").Value = "鯣鼷谭拼`t"
$ws.Range("C9").Value = "Q5encx`t"
$ws.Range("D9").Value = "2023-08-05 13:17:25`t"

$ws.Range("A10").Value = "8061`t"
$ws.Range("B10").Value = "塇鼄芋璆`t"
$ws.Range("C10").Value = "aUUI8S`t"
$ws.Range("D10").Value = "2023-08-05 13:17:25`t"

$ws.Range("A11").Value = "3282`t"
$ws.Range("B11").Value = "翲縬蒫湚`t"
$ws.Range("C11").Value = "6hp8Pq`t"
$ws.Range("D11").Value = "2023-08-05 13:17:25`t"

